$d = $word.ActiveDocument

$map = @{
  "79×71=5609" = "71×84=5964";
  "53×37=1961" = "42×88=3696";
  "54×42=2268" = "63×32=2016";
  "24×19=456"  = "60×36=2160";
  "58×73=4234" = "51×16=816";
  "94×92=8648" = "22×38=836";
  "98×85=8330" = "48×53=2544";
  "77×93=7161" = "41×71=2911";
  "88×86=7568" = "78×69=5382";
  "39×69=2691" = "87×38=3306";
  "32×63=2016" = "93×45=4185";
  "89×68=6052" = "42×38=1596";
  "74×67=4958" = "29×12=348";
  "21×91=1911" = "65×35=2275";
  "61×90=5490" = "77×73=5621";
  "15×87=1305" = "73×50=3650";
  "11×75=825"  = "23×48=1104";
  "90×13=1170" = "16×24=384";
  "54×95=5130" = "41×29=1189";
  "65×95=6175" = "75×14=1050";
  "12×80=960"  = "71×98=6958";
  "79×53=4187" = "22×63=1386";
  "79×63=4977" = "23×12=276";
  "62×83=5146" = "48×62=2976";
  "63×26=1638" = "24×80=1920";
}

foreach ($old in $map.Keys) {
  $new = $map[$old]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
